$wb = $excel.ActiveWorkbook

# The "imc" worksheet holds the Peso/Altura/IMC/Clasificación table.
$ws = $wb.Worksheets.Item("imc")

# Update the Peso (A) / Altura (B) inputs that changed; the IMC (C) and
# Clasificación (D) columns are driven by shared formulas and recompute
# automatically.
$ws.Range("B44").Value = 1
$ws.Range("B45").Value = 0.9
$ws.Range("A46").Value = 70
$ws.Range("B46").Value = 1.7
$ws.Range("B52").Value = 1.4
$ws.Range("B55").Value = 1.6
$ws.Range("A56").Value = 40
$ws.Range("B61").Value = 1
$ws.Range("A64").Value = 45
$ws.Range("B64").Value = 1.4
$ws.Range("B65").Value = 1.7
$ws.Range("B66").Value = 1.5
$ws.Range("A68").Value = 20
$ws.Range("B68").Value = 0.8
$ws.Range("B69").Value = 1.7
$ws.Range("B70").Value = 1.6
$ws.Range("A73").Value = 70
$ws.Range("B73").Value = 1.7
$ws.Range("B74").Value = 1.9
$ws.Range("A75").Value = 40
$ws.Range("A76").Value = 35
$ws.Range("B77").Value = 1.6
$ws.Range("A80").Value = 30
$ws.Range("A89").Value = 20
$ws.Range("B89").Value = 1
$ws.Range("B90").Value = 1.9
$ws.Range("B91").Value = 1.4
$ws.Range("A92").Value = 70
$ws.Range("A96").Value = 45
$ws.Range("A100").Value = 40

# Make "imc" the selected/active tab instead of "Hoja2".
$ws.Activate()
$ws.Select()
